$wb2 = $excel.ActiveWorkbook
$ws = $wb2.Worksheets.Item("primary")
# simple direct test: copy single cell C66 (numeric) to C200
$ws.Range("C66").Copy()
$ws.Range("C201").PasteSpecial()
Write-Output "done"
